# Auto-generated edit script: update market-price derived columns (H-N)
# across all profession sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1637.2
$ws.Range("I43").Value = 2270.889
$ws.Range("J43").Value = 686.6667
$ws.Range("K43").Value = 2270.889
$ws.Range("L43").Value = 686.6667
$ws.Range("M43").Value = -2201.889
$ws.Range("N43").Value = -824.6667
# Row 137
$ws.Range("H137").Value = 963.2222
$ws.Range("I137").Value = 806.5
$ws.Range("J137").Value = 1276.6666
$ws.Range("K137").Value = 2419.5
$ws.Range("L137").Value = 3829.9998
$ws.Range("M137").Value = 130.5
$ws.Range("N137").Value = -8929.9998
# Row 138
$ws.Range("H138").Value = 3095.88
$ws.Range("J138").Value = 3967.2876
$ws.Range("L138").Value = 11901.8628
$ws.Range("N138").Value = -22181.8628
# Row 141
$ws.Range("H141").Value = 3744.7
$ws.Range("I141").Value = 2761.25
$ws.Range("J141").Value = 4400.3335
$ws.Range("K141").Value = 8283.75
$ws.Range("L141").Value = 13201.0005
$ws.Range("M141").Value = -3103.75
$ws.Range("N141").Value = -23561.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2301.1304
$ws.Range("I2").Value = 2046.625
$ws.Range("J2").Value = 2882.8572
$ws.Range("K2").Value = 2046.625
$ws.Range("L2").Value = 2882.8572
$ws.Range("M2").Value = -1933.625
$ws.Range("N2").Value = -3108.8572
# Row 61
$ws.Range("H61").Value = 945.7646999999999
$ws.Range("I61").Value = 805.2
$ws.Range("K61").Value = 805.2
$ws.Range("M61").Value = -593.2
# Row 74
$ws.Range("H74").Value = 1003.3333
$ws.Range("I74").Value = 952.64
$ws.Range("J74").Value = 1161.75
$ws.Range("K74").Value = 952.64
$ws.Range("L74").Value = 1161.75
$ws.Range("M74").Value = -78.63999999999999
$ws.Range("N74").Value = -2909.75
# Row 77
$ws.Range("H77").Value = 1003.3333
$ws.Range("I77").Value = 952.64
$ws.Range("J77").Value = 1161.75
$ws.Range("K77").Value = 4763.2
$ws.Range("L77").Value = 5808.75
$ws.Range("M77").Value = -395.1999999999998
$ws.Range("N77").Value = -14544.75
# Row 110
$ws.Range("H110").Value = 629.08
$ws.Range("I110").Value = 613.625
$ws.Range("K110").Value = 613.625
$ws.Range("M110").Value = 1431.375
# Row 116
$ws.Range("H116").Value = 2301.1304
$ws.Range("I116").Value = 2046.625
$ws.Range("J116").Value = 2882.8572
$ws.Range("K116").Value = 2046.625
$ws.Range("L116").Value = 2882.8572
$ws.Range("M116").Value = 247.375
$ws.Range("N116").Value = -7470.8572
# Row 122
$ws.Range("H122").Value = 985.3333
$ws.Range("I122").Value = 978
$ws.Range("K122").Value = 2934
$ws.Range("M122").Value = -484
# Row 132
$ws.Range("H132").Value = 2465.7856
$ws.Range("I132").Value = 1384.8
$ws.Range("J132").Value = 3066.3333
$ws.Range("K132").Value = 4154.4
$ws.Range("L132").Value = 9198.999899999999
$ws.Range("M132").Value = -1624.4
$ws.Range("N132").Value = -14258.9999
# Row 136
$ws.Range("H136").Value = 945.7646999999999
$ws.Range("I136").Value = 805.2
$ws.Range("K136").Value = 2415.6
$ws.Range("M136").Value = 134.3999999999996

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2301.1304
$ws.Range("I3").Value = 2046.625
$ws.Range("J3").Value = 2882.8572
$ws.Range("K3").Value = 2046.625
$ws.Range("L3").Value = 2882.8572
$ws.Range("M3").Value = -1932.625
$ws.Range("N3").Value = -3110.8572
# Row 134
$ws.Range("H134").Value = 67692.25999999999
$ws.Range("I134").Value = 3293.3333
$ws.Range("J134").Value = 202930
$ws.Range("K134").Value = 9879.999899999999
$ws.Range("L134").Value = 608790
$ws.Range("M134").Value = -7344.999899999999
$ws.Range("N134").Value = -613860

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 1390
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1390
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1390
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -1614
# Row 31
$ws.Range("H31").Value = 3706.0386
$ws.Range("I31").Value = 3845.348
$ws.Range("K31").Value = 3845.348
$ws.Range("M31").Value = -3550.348
# Row 34
$ws.Range("H34").Value = 3706.0386
$ws.Range("I34").Value = 3845.348
$ws.Range("K34").Value = 3845.348
$ws.Range("M34").Value = -3643.348
# Row 58
$ws.Range("H58").Value = 6253.2
$ws.Range("I58").Value = 1212.6154
$ws.Range("J58").Value = 15614.286
$ws.Range("K58").Value = 1212.6154
$ws.Range("L58").Value = 15614.286
$ws.Range("M58").Value = -1009.6154
$ws.Range("N58").Value = -16020.286
# Row 132
$ws.Range("H132").Value = 2025.25
$ws.Range("I132").Value = 1060.7778
$ws.Range("J132").Value = 2814.3635
$ws.Range("K132").Value = 3182.3334
$ws.Range("L132").Value = 8443.0905
$ws.Range("M132").Value = -652.3334000000004
$ws.Range("N132").Value = -13503.0905
# Row 134
$ws.Range("H134").Value = 2777.0322
$ws.Range("I134").Value = 2360.818
$ws.Range("J134").Value = 3794.4443
$ws.Range("K134").Value = 7082.454000000001
$ws.Range("L134").Value = 11383.3329
$ws.Range("M134").Value = -4547.454000000001
$ws.Range("N134").Value = -16453.3329
# Row 136
$ws.Range("H136").Value = 6253.2
$ws.Range("I136").Value = 1212.6154
$ws.Range("J136").Value = 15614.286
$ws.Range("K136").Value = 3637.8462
$ws.Range("L136").Value = 46842.858
$ws.Range("M136").Value = -1087.8462
$ws.Range("N136").Value = -51942.858

$ws = $wb.Worksheets.Item("CUL")
# Row 104
$ws.Range("H104").Value = 5482.5
$ws.Range("J104").Value = 5482.5
$ws.Range("L104").Value = 16447.5
$ws.Range("N104").Value = -21689.5
# Row 113
$ws.Range("H113").Value = 747.7059
$ws.Range("J113").Value = 776.23334
$ws.Range("L113").Value = 2328.70002
$ws.Range("N113").Value = -6668.70002

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5898.4443
$ws.Range("J70").Value = 7000
$ws.Range("L70").Value = 7000
$ws.Range("N70").Value = -7540
# Row 73
$ws.Range("H73").Value = 5898.4443
$ws.Range("J73").Value = 7000
$ws.Range("L73").Value = 7000
$ws.Range("N73").Value = -8872
# Row 132
$ws.Range("H132").Value = 3109.111
$ws.Range("I132").Value = 3008.7058
$ws.Range("J132").Value = 3279.8
$ws.Range("K132").Value = 9026.117400000001
$ws.Range("L132").Value = 9839.400000000001
$ws.Range("M132").Value = -6496.117400000001
$ws.Range("N132").Value = -14899.4

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4846.077
$ws.Range("I122").Value = 7364.8125
$ws.Range("J122").Value = 3093.913
$ws.Range("K122").Value = 22094.4375
$ws.Range("L122").Value = 9281.739
$ws.Range("M122").Value = -19644.4375
$ws.Range("N122").Value = -14181.739
# Row 132
$ws.Range("H132").Value = 2272.9546
$ws.Range("I132").Value = 1392.963
$ws.Range("K132").Value = 4178.889
$ws.Range("M132").Value = -1648.889
# Row 136
$ws.Range("H136").Value = 4666.125
$ws.Range("I136").Value = 1443.7222
$ws.Range("J136").Value = 14333.333
$ws.Range("K136").Value = 4331.1666
$ws.Range("L136").Value = 42999.999
$ws.Range("M136").Value = -1781.1666
$ws.Range("N136").Value = -48099.999

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 502020
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 502020
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 502020
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -502246
# Row 132
$ws.Range("H132").Value = 1134.4828
$ws.Range("I132").Value = 1072.6923
$ws.Range("J132").Value = 1670
$ws.Range("K132").Value = 3218.0769
$ws.Range("L132").Value = 5010
$ws.Range("M132").Value = -688.0769
$ws.Range("N132").Value = -10070
# Row 136
$ws.Range("H136").Value = 2172.8
$ws.Range("I136").Value = 2333.8572
$ws.Range("J136").Value = 1528.5714
$ws.Range("K136").Value = 7001.571599999999
$ws.Range("L136").Value = 4585.7142
$ws.Range("M136").Value = -4451.571599999999
$ws.Range("N136").Value = -9685.7142
